$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the oldest data row (the row that used to start the series at A2=39400).
# This shifts every subsequent row up by one and keeps the date-format style intact.
$ws.Rows(2).Delete()

# After the shift, column E (y_1_forecast) needs to be cleared for the first four
# data rows (the forecaster no longer has a value there) and recomputed for the rest.
$ws.Range("E2:E5").ClearContents()

# Refresh column C (y_0_forecast) value on row 3 with the newly recomputed figure.
$ws.Range("C3").Value = -1.317619545389281

# Recomputed column E (y_1_forecast) values for rows 6-18.
$ws.Range("E6").Value = 4.792854588620821
$ws.Range("E7").Value = 7.239454936865775
$ws.Range("E8").Value = 0.4361429468412448
$ws.Range("E9").Value = 1.392195163617171
$ws.Range("E10").Value = 1.493220091771108
$ws.Range("E11").Value = 2.622364272988187
$ws.Range("E12").Value = 3.113086948791377
$ws.Range("E13").Value = 2.757421718286168
$ws.Range("E14").Value = -1.373617952268746
$ws.Range("E15").Value = 1.324283050325015
$ws.Range("E16").Value = -0.5025420863900898
$ws.Range("E17").Value = -2.371854438773213
$ws.Range("E18").Value = -1.436963918858969
